$wb = $excel.ActiveWorkbook
$cols = @("B", "C", "D", "E", "F", "G", "H")

# --- Sheet "900-1145 AM" : constant row pattern 9,10,11,12,13,14,15 for B:H, rows 2-13 ---
$ws = $wb.Worksheets.Item("900-1145 AM")
$vals = @(9, 10, 11, 12, 13, 14, 15)
for ($r = 2; $r -le 13; $r++) {
    for ($ci = 0; $ci -lt $cols.Length; $ci++) {
        $ws.Range("$($cols[$ci])$r").Value = $vals[$ci]
    }
}
$ws.Range("H2:H13").Select()

# --- Sheet "1200-245 PM" : constant value per row across B:H, increasing by 1 each row (12..23) ---
$ws = $wb.Worksheets.Item("1200-245 PM")
for ($r = 2; $r -le 13; $r++) {
    $val = $r + 10
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $val
    }
}
$ws.Range("B13:H13").Select()

# --- Sheet "300-545 PM" : constant value 100 across all of B2:H13 ---
$ws = $wb.Worksheets.Item("300-545 PM")
for ($r = 2; $r -le 13; $r++) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = 100
    }
}
$ws.Range("G16").Select()

# --- Sheet "600-845 PM" : diagonal pattern value = (row-2) + (colIndex-1), with a few manual gaps/typos ---
$ws = $wb.Worksheets.Item("600-845 PM")
for ($r = 2; $r -le 13; $r++) {
    for ($ci = 0; $ci -lt $cols.Length; $ci++) {
        $col = $cols[$ci]
        $val = ($r - 2) + ($ci + 1)
        $ws.Range("$col$r").Value = $val
    }
}
# manual omissions / typos present in the authored data
$ws.Range("C9").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("D11").Value = 16
$ws.Range("E11").Value = " "
$ws.Range("F12").ClearContents()
$ws.Range("G13").ClearContents()

$ws.Activate()
$ws.Range("H14").Select()
